$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing content/formatting across the previous used range (A1:G65)
$ws.Range("A1:G65").Clear()

# --- Plain (default-style) cell values ---
$ws.Range("A1").Value = "VNFD"
$ws.Range("B11").Value = "id"
$ws.Range("B12").Value = "name"
$ws.Range("B13").Value = "description"
$ws.Range("B14").Value = "int-cpd []"
$ws.Range("D15").Value = "id"
$ws.Range("D17").Value = "layer-protocol"
$ws.Range("E18").Value = "x"
$ws.Range("G19").Value = "x"
$ws.Range("G20").Value = "x"
$ws.Range("B22").Value = "boot-order"
$ws.Range("E22").Value = "x"
$ws.Range("D23").Value = "value"
$ws.Range("C24").Value = "x"
$ws.Range("C25").Value = "x"
$ws.Range("C26").Value = "x"
$ws.Range("C27").Value = "x"
$ws.Range("D27").Value = "Missing All node of this list"
$ws.Range("C28").Value = "x"
$ws.Range("A30").Value = "sw-image-desc []"
$ws.Range("B30").Value = "id"
$ws.Range("D30").Value = "value need to change"
$ws.Range("B31").Value = "name"
$ws.Range("D31").Value = "value need to change"
$ws.Range("B32").Value = "image-name-variable"
$ws.Range("D32").Value = "we can always populate this field"
$ws.Range("B33").Value = "version"
$ws.Range("B34").Value = "checksum"
$ws.Range("B35").Value = "container-format"
$ws.Range("B36").Value = "disk-format"
$ws.Range("B37").Value = "min-disk"
$ws.Range("B38").Value = "size"
$ws.Range("B39").Value = "image"
$ws.Range("A40").Value = "virtual-compute-descriptor []"
$ws.Range("B41").Value = "id"
$ws.Range("D41").Value = "value might need to change"
$ws.Range("B42").Value = "flavor-name-variable"
$ws.Range("B43").Value = "virtual-cpu"
$ws.Range("B44").Value = "virtual-memory"
$ws.Range("A45").Value = "virtual-storage-descriptor []"
$ws.Range("B46").Value = "id"
$ws.Range("B47").Value = "type-of-storage"
$ws.Range("B48").Value = "size-of-storage"
$ws.Range("B49").Value = "sw-image-desc"
$ws.Range("A50").Value = "int-virtual-link-desc[]"
$ws.Range("A52").Value = "configurable-properties"
$ws.Range("A53").Value = "df"
$ws.Range("B54").Value = "id"
$ws.Range("B55").Value = "description"
$ws.Range("B56").Value = "vdu-profile []"
$ws.Range("D57").Value = "id"
$ws.Range("D58").Value = "min-number-of-instances"
$ws.Range("D59").Value = "max-number-of-instances"
$ws.Range("E60").Value = "x"
$ws.Range("C61").Value = "x"
$ws.Range("D61").Value = "Missing All node of this list"
$ws.Range("C62").Value = "x"
$ws.Range("D62").Value = "Missing All node of this list"
$ws.Range("C63").Value = "x"
$ws.Range("D63").Value = "Missing All node of this list"

# --- Cells needing the red "missing/TBD" style (matches style index 1: red font) ---
$ws.Range("D18").Value = "additional-sol1-parameters"
$ws.Range("D18").Font.Color = 255
$ws.Range("F19").Value = "allowed-address-variable"
$ws.Range("F19").Font.Color = 255
$ws.Range("F20").Value = "security-group-variable"
$ws.Range("F20").Font.Color = 255
$ws.Range("D22").Value = "key"
$ws.Range("D22").Font.Color = 255
$ws.Range("B24").Value = "virtual-compute-desc"
$ws.Range("B24").Font.Color = 255
$ws.Range("B25").Value = "virtual-storage-desc"
$ws.Range("B25").Font.Color = 255
$ws.Range("B26").Value = "sw-image-desc"
$ws.Range("B26").Font.Color = 255
$ws.Range("B27").Value = "day0 []"
$ws.Range("B27").Font.Color = 255
$ws.Range("B28").Value = "device-type"
$ws.Range("B28").Font.Color = 255
$ws.Range("D60").Value = "affinity-or-anti-affinity-group []"
$ws.Range("D60").Font.Color = 255
$ws.Range("B61").Value = "instantiation-level []"
$ws.Range("B61").Font.Color = 255
$ws.Range("B62").Value = "scaling-aspect []"
$ws.Range("B62").Font.Color = 255
$ws.Range("B63").Value = "affinity-or-anti-affinity-group []"
$ws.Range("B63").Font.Color = 255

# --- Cells needing the header-like (style 2) formatting ---
# Seed cell A2 already carries style-2 formatting from the original template;
# copy its format onto the other cells that need the same look.
$ws.Range("A2").Value = "id"
$ws.Range("A3").Value = "provider"
$ws.Range("A4").Value = "product-name"
$ws.Range("A5").Value = "software-version"
$ws.Range("A6").Value = "version"
$ws.Range("A7").Value = "product-info-name"
$ws.Range("A8").Value = "product-info-description"
$ws.Range("A9").Value = "vnfm-info"
$ws.Range("A10").Value = "vdu []"
$ws.Range("D16").Value = "int-virtual-link-desc"
$ws.Range("A51").Value = "ext-cpd []"
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("A51").PasteSpecial(-4122)

# --- Sheet view / selection state ---
$ws.Range("B61").Select()
